# ListingPageTestData.xlsx - "Test cases added as suggested"
# Adds two new sheets (article, filter) with new test-case rows, tweaks a
# couple of existing testcase descriptions, and nudges a few view-state
# bits (selection, active sheet).

$wb = $excel.ActiveWorkbook

$xlLeft = -4131
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. bannerTexts: fix up two expected-text descriptions + row height
# ---------------------------------------------------------------------
$bannerTexts = $wb.Worksheets.Item("bannerTexts")
$bannerTexts.Range("A2").Value = "Verify the banner heading text"
$bannerTexts.Range("A3").Value = "Verify the banner sub-heading text"
$bannerTexts.Rows.Item(3).RowHeight = 30

# ---------------------------------------------------------------------
# 2. quickFilter: widen column A a touch
# ---------------------------------------------------------------------
$quickFilter = $wb.Worksheets.Item("quickFilter")
$quickFilter.Columns.Item(1).ColumnWidth = 33.65

# ---------------------------------------------------------------------
# 3. Create the two new sheets (added right after "sortBy").
#    Both are inserted immediately after "sortBy" - doing "filter" first
#    means "article" ends up between "sortBy" and "filter" in tab order,
#    matching the target layout: ... sortBy, article, filter.
# ---------------------------------------------------------------------
$sortBy = $wb.Worksheets.Item("sortBy")

$newSheet1 = $wb.Worksheets.Add($null, $sortBy)
$newSheet1.Name = "filter"

$newSheet2 = $wb.Worksheets.Add($null, $sortBy)
$newSheet2.Name = "article"

# Sheet references above track tab *position*, not identity, so once both
# inserts are done, re-resolve them by name before touching their content.
$filter = $wb.Worksheets.Item("filter")
$article = $wb.Worksheets.Item("article")

# ---------------------------------------------------------------------
# 4. "filter" sheet content
# ---------------------------------------------------------------------
$filter.Columns.Item(1).ColumnWidth = 48.3047619047619
$filter.Columns.Item(2).ColumnWidth = 9.16666666666667
$filter.Columns.Item(3).ColumnWidth = 34.3047619047619
$filter.Columns.Item(4).ColumnWidth = 23.7380952380952

$filter.Range("A1").Value = "testcase"
$filter.Range("A1").HorizontalAlignment = $xlLeft
$filter.Range("A1").VerticalAlignment = $xlCenter
$filter.Range("B1").Value = "execution"
$filter.Range("B1").HorizontalAlignment = $xlLeft
$filter.Range("B1").VerticalAlignment = $xlCenter
$filter.Range("C1").Value = "filterText"
$filter.Range("C1").HorizontalAlignment = $xlLeft
$filter.Range("C1").VerticalAlignment = $xlCenter
$filter.Range("D1").Value = "subFilter"

$filterRows = @(
    @("Verify Ages filter for 2+", "yes", "Ages", "2+"),
    @("Verify Ages filter for 6+", "yes", "Ages", "6+"),
    @("Verify Ages filter for 6-11", "yes", "Ages", "6-11"),
    @("Verify Ages filter 12+", "yes", "Ages", "12+"),
    @("Verify Symptom filter for Itchy Nose and/or Throat", "yes", "Symptom", "Itchy Nose and/or Throat"),
    @("Verify Symptom filter for Itchy Skin", "yes", "Symptom", "Itchy Skin"),
    @("Verify Symptom filter for Itchy, Watery Eyes", "yes", "Symptom", "Itchy, Watery Eyes"),
    @("Verify Symptom filter for Nasal Congestion", "yes", "Symptom", "Nasal Congestion"),
    @("Verify Symptom filter for Runny Nose", "yes", "Symptom", "Runny Nose"),
    @("Verify Symptom filter for Scrapes & Minor Cuts", "yes", "Symptom", "Scrapes & Minor Cuts"),
    @("Verify Symptom filter for Sinus Pressure", "yes", "Symptom", "Sinus Pressure"),
    @("Verify Symptom filter for Sneezing", "yes", "Symptom", "Sneezing"),
    @("Verify Symptom filter for Sunburn & Minor Burns", "yes", "Symptom", "Sunburn & Minor Burns"),
    @("Verify Category filter for Adult Products", "yes", "Category", "Adult Products"),
    @("Verify Category filter for Children's Products", "yes", "Category", "Children's Products"),
    @("Verify Category filter for Topical Products", "yes", "Category", "Topical Products")
)

$r = 2
foreach ($row in $filterRows) {
    $filter.Cells.Item($r, 1).Value = $row[0]
    $filter.Cells.Item($r, 2).Value = $row[1]
    $filter.Cells.Item($r, 3).Value = $row[2]
    $filter.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# "6-11" looks like a date - force the cell to plain text like the
# source file does, so it round-trips as the string "6-11".
$filter.Range("D4").NumberFormat = "@"
$filter.Range("D4").Value = "6-11"

# ---------------------------------------------------------------------
# 5. "article" sheet content
# ---------------------------------------------------------------------
$article.Columns.Item(1).ColumnWidth = 28.7380952380952
$article.Columns.Item(2).ColumnWidth = 8.73809523809524
$article.Columns.Item(3).ColumnWidth = 46.8571428571429

$article.Range("A1").Value = "testcase"
$article.Range("A1").HorizontalAlignment = $xlLeft
$article.Range("A1").VerticalAlignment = $xlCenter
$article.Range("B1").Value = "execution"
$article.Range("B1").HorizontalAlignment = $xlLeft
$article.Range("B1").VerticalAlignment = $xlCenter
$article.Range("C1").Value = "altTxt"
$article.Range("D1").Value = "url"
$article.Range("D1").HorizontalAlignment = $xlLeft
$article.Range("D1").VerticalAlignment = $xlCenter

$article.Range("A2").Value = "Verify dosing guide article"
$article.Range("B2").Value = "yes"
$article.Range("C2").Value = "Benadryl dosing for children and adults"
$article.Hyperlinks.Add($article.Range("D2"), "https://www.benadryl.com/benadryl-dosing-guide", [System.Type]::Missing, [System.Type]::Missing, "https://www.benadryl.com/benadryl-dosing-guide")

$article.Range("A3").Value = "Verify cold and allergies article"
$article.Range("B3").Value = "yes"
$article.Range("C3").Value = " Woman with cold or allergies"
$article.Hyperlinks.Add($article.Range("D3"), "https://www.benadryl.com/cold/cold-or-allergies", [System.Type]::Missing, [System.Type]::Missing, "https://www.benadryl.com/cold/cold-or-allergies")

# ---------------------------------------------------------------------
# 6. View-state touch ups on the original sheets
#    (re-resolve every sheet by name - see the note above about stale
#    position-based references)
# ---------------------------------------------------------------------
$wb.Worksheets.Item("bannerImg").Range("B27").Select()

$wb.Worksheets.Item("bannerTexts").Range("A8:A9").Select()

$wb.Worksheets.Item("quickFilter").Range("C31").Select()

$wb.Worksheets.Item("sortBy").Range("E21").Select()

$wb.Worksheets.Item("article").Range("C24").Select()

# Make "filter" the active tab, matching the saved workbook view.
$finalFilter = $wb.Worksheets.Item("filter")
$finalFilter.Activate()
$finalFilter.Range("C2").Select()
